# Updated cryptos list on Thu Jul  4 11:55:49 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

Set-TextCell 2 4 "57.796.42"
Set-TextCell 2 5 "  -4.15%  "
Set-TextCell 3 4 "3.159.27"
Set-TextCell 3 5 "  -4.39%  "
Set-TextCell 4 5 "  +0.06%  "
Set-TextCell 5 4 "525.57"
Set-TextCell 5 5 "  -6.18%  "
Set-TextCell 6 4 "133.18"
Set-TextCell 6 5 "  -7.37%  "
Set-TextCell 7 5 "  -0.03%  "
Set-TextCell 8 4 "3.156.19"
Set-TextCell 8 5 "  -4.57%  "
Set-TextCell 9 5 "  -6.10%  "
Set-TextCell 10 4 "7.28"
Set-TextCell 10 5 "  -7.22%  "
Set-TextCell 11 5 "  -6.88%  "
Set-TextCell 12 5 "  -4.23%  "
Set-TextCell 13 4 "3.698.46"
Set-TextCell 13 5 "  -4.57%  "
Set-TextCell 14 5 "  -1.77%  "
Set-TextCell 15 4 "25.92"
Set-TextCell 15 5 "  -5.13%  "
Set-TextCell 16 4 "3.158.01"
Set-TextCell 16 5 "  -4.30%  "
Set-TextCell 17 4 "57.788.32"
Set-TextCell 17 5 "  -4.17%  "
Set-TextCell 18 5 "  -7.94%  "
Set-TextCell 19 4 "5.81"
Set-TextCell 19 5 "  -5.05%  "
Set-TextCell 20 4 "13.05"
Set-TextCell 20 5 "  -9.22%  "
Set-TextCell 21 4 "8.07"
Set-TextCell 21 5 "  -6.21%  "
Set-TextCell 22 4 "346.18"
Set-TextCell 22 5 "  -7.51%  "
Set-TextCell 23 5 "  -0.23%  "
Set-TextCell 24 4 "69.73"
Set-TextCell 24 5 "  -5.43%  "
Set-TextCell 25 5 "  -6.68%  "
Set-TextCell 26 4 "3.285.37"
Set-TextCell 26 5 "  -5.05%  "
Set-TextCell 27 4 "0.0₃0960"
Set-TextCell 27 5 "  -8.55%  "
Set-TextCell 28 5 "  -4.02%  "
Set-TextCell 29 4 "1.00"
Set-TextCell 29 5 "  +0.26%  "
Set-TextCell 30 4 "6.87"
Set-TextCell 30 5 "  -5.02%  "
Set-TextCell 31 4 "0.998"
Set-TextCell 31 5 "  -0.20%  "
Set-TextCell 32 5 "  -8.21%  "
Set-TextCell 33 4 "6.93"
Set-TextCell 33 5 "  -9.29%  "
Set-TextCell 34 5 "  -3.82%  "
Set-TextCell 35 5 "  -5.33%  "
Set-TextCell 36 4 "4.91"
Set-TextCell 36 5 "  -5.38%  "
Set-TextCell 37 4 "159.36"
Set-TextCell 37 5 "  -4.49%  "
Set-TextCell 38 5 "  -7.48%  "
Set-TextCell 39 5 "  -7.74%  "
Set-TextCell 40 4 "26.01"
Set-TextCell 40 5 "  -5.35%  "
Set-TextCell 41 5 "  -5.39%  "
Set-TextCell 42 4 "3.185.01"
Set-TextCell 42 5 "  -4.58%  "
Set-TextCell 43 4 "40.68"
Set-TextCell 43 5 "  -2.79%  "
Set-TextCell 44 4 "0.699"
Set-TextCell 44 5 "  -7.12%  "
Set-TextCell 45 5 "  -4.15%  "
Set-TextCell 46 4 "3.95"
Set-TextCell 46 5 "  -6.13%  "
Set-TextCell 47 4 "0.999"
Set-TextCell 47 5 "  -0.04%  "
Set-TextCell 48 4 "1.47"
Set-TextCell 48 5 "  -8.35%  "
Set-TextCell 49 4 "2.269.57"
Set-TextCell 49 5 "  -4.30%  "
Set-TextCell 50 5 "  -5.99%  "
Set-TextCell 51 4 "20.51"
